$d = $word.ActiveDocument

# Update the date heading (first paragraph)
$d.Paragraphs.Item(1).Range.Text = "2025-12-14 Sunday"

# Update each math-problem cell in the table, addressed positionally
# to correctly disambiguate duplicate cell text values.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "32+49="
$t.Cell(1, 2).Range.Text = "21-8="
$t.Cell(1, 3).Range.Text = "16+27="
$t.Cell(1, 4).Range.Text = "18+73="
$t.Cell(1, 5).Range.Text = "91-83="
$t.Cell(2, 1).Range.Text = "88+5="
$t.Cell(2, 2).Range.Text = "18+66="
$t.Cell(2, 3).Range.Text = "90-7="
$t.Cell(2, 4).Range.Text = "83-76="
$t.Cell(2, 5).Range.Text = "13+19="
$t.Cell(3, 1).Range.Text = "17+6="
$t.Cell(3, 2).Range.Text = "41-2="
$t.Cell(3, 3).Range.Text = "28+34="
$t.Cell(3, 4).Range.Text = "19+15="
$t.Cell(3, 5).Range.Text = "58+4="
$t.Cell(4, 1).Range.Text = "8+13="
$t.Cell(4, 2).Range.Text = "57+26="
$t.Cell(4, 3).Range.Text = "52-28="
$t.Cell(4, 4).Range.Text = "26+69="
$t.Cell(4, 5).Range.Text = "46+46="
$t.Cell(5, 1).Range.Text = "48+48="
$t.Cell(5, 2).Range.Text = "8+14="
$t.Cell(5, 3).Range.Text = "71-34="
$t.Cell(5, 4).Range.Text = "42-7="
$t.Cell(5, 5).Range.Text = "34-19="
$t.Cell(6, 1).Range.Text = "95-89="
$t.Cell(6, 2).Range.Text = "83-18="
$t.Cell(6, 3).Range.Text = "17+19="
$t.Cell(6, 4).Range.Text = "91-72="
$t.Cell(6, 5).Range.Text = "15+66="
$t.Cell(7, 1).Range.Text = "34+49="
$t.Cell(7, 2).Range.Text = "70-37="
$t.Cell(7, 3).Range.Text = "58+8="
$t.Cell(7, 4).Range.Text = "35-17="
$t.Cell(7, 5).Range.Text = "8+24="
$t.Cell(8, 1).Range.Text = "88-79="
$t.Cell(8, 2).Range.Text = "81-7="
$t.Cell(8, 3).Range.Text = "81-75="
$t.Cell(8, 4).Range.Text = "52-33="
$t.Cell(8, 5).Range.Text = "77+14="
$t.Cell(9, 1).Range.Text = "82-65="
$t.Cell(9, 2).Range.Text = "66-28="
$t.Cell(9, 3).Range.Text = "62-17="
$t.Cell(9, 4).Range.Text = "9+26="
$t.Cell(9, 5).Range.Text = "69+7="
$t.Cell(10, 1).Range.Text = "19+25="
$t.Cell(10, 2).Range.Text = "76+18="
$t.Cell(10, 3).Range.Text = "45+38="
$t.Cell(10, 4).Range.Text = "38+59="
$t.Cell(10, 5).Range.Text = "27+57="
$t.Cell(11, 1).Range.Text = "6+7="
$t.Cell(11, 2).Range.Text = "27+68="
$t.Cell(11, 3).Range.Text = "25+36="
$t.Cell(11, 4).Range.Text = "17+45="
$t.Cell(11, 5).Range.Text = "90-66="
$t.Cell(12, 1).Range.Text = "87+6="
$t.Cell(12, 2).Range.Text = "29+64="
$t.Cell(12, 3).Range.Text = "51-12="
$t.Cell(12, 4).Range.Text = "29+15="
$t.Cell(12, 5).Range.Text = "82-14="
$t.Cell(13, 1).Range.Text = "28+44="
$t.Cell(13, 2).Range.Text = "15+58="
$t.Cell(13, 3).Range.Text = "13+49="
$t.Cell(13, 4).Range.Text = "28+29="
$t.Cell(13, 5).Range.Text = "90-79="
$t.Cell(14, 1).Range.Text = "16+77="
$t.Cell(14, 2).Range.Text = "82-46="
$t.Cell(14, 3).Range.Text = "70-34="
$t.Cell(14, 4).Range.Text = "98-59="
$t.Cell(14, 5).Range.Text = "54-38="
$t.Cell(15, 1).Range.Text = "6+48="
$t.Cell(15, 2).Range.Text = "37+18="
$t.Cell(15, 3).Range.Text = "53-44="
$t.Cell(15, 4).Range.Text = "27+46="
$t.Cell(15, 5).Range.Text = "20-18="
$t.Cell(16, 1).Range.Text = "42-38="
$t.Cell(16, 2).Range.Text = "64+9="
$t.Cell(16, 3).Range.Text = "75+8="
$t.Cell(16, 4).Range.Text = "57+7="
$t.Cell(16, 5).Range.Text = "47+37="
$t.Cell(17, 1).Range.Text = "90-47="
$t.Cell(17, 2).Range.Text = "96-49="
$t.Cell(17, 3).Range.Text = "85+8="
$t.Cell(17, 4).Range.Text = "8+49="
$t.Cell(17, 5).Range.Text = "9+25="
$t.Cell(18, 1).Range.Text = "25+28="
$t.Cell(18, 2).Range.Text = "22-15="
$t.Cell(18, 3).Range.Text = "97-9="
$t.Cell(18, 4).Range.Text = "80-54="
$t.Cell(18, 5).Range.Text = "81-24="
$t.Cell(19, 1).Range.Text = "15+48="
$t.Cell(19, 2).Range.Text = "19+33="
$t.Cell(19, 3).Range.Text = "19+33="
$t.Cell(19, 4).Range.Text = "14+18="
$t.Cell(19, 5).Range.Text = "4+17="
$t.Cell(20, 1).Range.Text = "86-69="
$t.Cell(20, 2).Range.Text = "40-24="
$t.Cell(20, 3).Range.Text = "84-56="
$t.Cell(20, 4).Range.Text = "56-17="
$t.Cell(20, 5).Range.Text = "74-55="
